# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E30) listed account-statement periods in
# descending order (2001 .. 1811). This update refreshes/reorders the data
# so the periods run in ascending chronological order (1811 .. 2001), as
# part of adding the first batch of the new account-statement periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$periodosMora = @(
    "1811", "1812",
    "1901", "1902", "1903", "1904", "1905", "1906",
    "1907", "1908", "1909", "1910", "1911", "1912",
    "2001"
)

$firstRow = 16
for ($i = 0; $i -lt $periodosMora.Length; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $periodosMora[$i]
}
